$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 1: shift labels C1/D1/E1 to the right (wrapping)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update row 2 values accordingly
$ws.Range("C2").Value = "s__Amedibacillus dolichus"
$ws.Range("D2").Value = "s__Amedibacillus dolichus"
$ws.Range("E2").Value = 1
